$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Cells.Item(3, 6).Value = 2.96
$ws.Cells.Item(3, 7).Value = 3.3
$ws.Cells.Item(3, 8).Value = 2.4
$ws.Cells.Item(3, 9).Value = 2.6
$ws.Cells.Item(3, 10).Value = 3.45
$ws.Cells.Item(3, 11).Value = 3.8

# Row 4
$ws.Cells.Item(4, 6).Value = 4.3
$ws.Cells.Item(4, 7).Value = 5.1
$ws.Cells.Item(4, 8).Value = 2.08
$ws.Cells.Item(4, 9).Value = 2.26
$ws.Cells.Item(4, 10).Value = 2.96
$ws.Cells.Item(4, 11).Value = 3.25
$ws.Cells.Item(4, 14).Value = 2.5
$ws.Cells.Item(4, 15).Value = 1.54
$ws.Cells.Item(4, 16).Value = 1.51
$ws.Cells.Item(4, 17).Value = 2.64
$ws.Cells.Item(4, 20).Value = 2.12
$ws.Cells.Item(4, 21).Value = 1.71
$ws.Cells.Item(4, 24).Value = 8.800000000000001
$ws.Cells.Item(4, 25).Value = 6.8
$ws.Cells.Item(4, 26).Value = 14
$ws.Cells.Item(4, 27).Value = 32
$ws.Cells.Item(4, 28).Value = 13.5
$ws.Cells.Item(4, 29).Value = 7.6
$ws.Cells.Item(4, 31).Value = 38
$ws.Cells.Item(4, 32).Value = 34
$ws.Cells.Item(4, 33).Value = 22
$ws.Cells.Item(4, 34).Value = 32
$ws.Cells.Item(4, 35).Value = 65
$ws.Cells.Item(4, 36).Value = 140
$ws.Cells.Item(4, 37).Value = 85
$ws.Cells.Item(4, 38).Value = 130
$ws.Cells.Item(4, 39).Value = 250
$ws.Cells.Item(4, 40).Value = 150
$ws.Cells.Item(4, 41).Value = 32

# Row 5
$ws.Cells.Item(5, 6).Value = 2.92
$ws.Cells.Item(5, 7).Value = 3.45
$ws.Cells.Item(5, 8).Value = 2.66
$ws.Cells.Item(5, 9).Value = 3.05
$ws.Cells.Item(5, 10).Value = 2.9
$ws.Cells.Item(5, 11).Value = 3.3
$ws.Cells.Item(5, 16).Value = 1.55

# Row 6
$ws.Cells.Item(6, 9).Value = 2.52
$ws.Cells.Item(6, 10).Value = 3.15
$ws.Cells.Item(6, 16).Value = 1.76
$ws.Cells.Item(6, 17).Value = 2.06

# Row 8
$ws.Cells.Item(8, 8).Value = 1.67
$ws.Cells.Item(8, 9).Value = 1.81
$ws.Cells.Item(8, 16).Value = 1.67
$ws.Cells.Item(8, 17).Value = 2.3

# Row 9
$ws.Cells.Item(9, 9).Value = 2.18

# Row 10
$ws.Cells.Item(10, 9).Value = 2.26
$ws.Cells.Item(10, 16).Value = 1.86

# Row 11
$ws.Cells.Item(11, 7).Value = 4.7
$ws.Cells.Item(11, 8).Value = 1.92
$ws.Cells.Item(11, 9).Value = 1.98
$ws.Cells.Item(11, 10).Value = 3.5
$ws.Cells.Item(11, 21).Value = 2.02
$ws.Cells.Item(11, 24).Value = 15.5
$ws.Cells.Item(11, 26).Value = 13
$ws.Cells.Item(11, 27).Value = 26
$ws.Cells.Item(11, 29).Value = 9.199999999999999
$ws.Cells.Item(11, 34).Value = 22
$ws.Cells.Item(11, 35).Value = 42
$ws.Cells.Item(11, 40).Value = 85
$ws.Cells.Item(11, 41).Value = 16

# Row 12
$ws.Cells.Item(12, 6).Value = 1.32
$ws.Cells.Item(12, 7).Value = 1.37
$ws.Cells.Item(12, 8).Value = 11
$ws.Cells.Item(12, 9).Value = 12.5
$ws.Cells.Item(12, 10).Value = 5.6
$ws.Cells.Item(12, 11).Value = 6.4
$ws.Cells.Item(12, 16).Value = 2.46
$ws.Cells.Item(12, 20).Value = 2.02
$ws.Cells.Item(12, 21).Value = 1.85
$ws.Cells.Item(12, 27).Value = 470
$ws.Cells.Item(12, 28).Value = 9.800000000000001
$ws.Cells.Item(12, 31).Value = 200
$ws.Cells.Item(12, 32).Value = 8.800000000000001
$ws.Cells.Item(12, 35).Value = 150
$ws.Cells.Item(12, 36).Value = 11
$ws.Cells.Item(12, 39).Value = 160

# Row 13
$ws.Cells.Item(13, 10).Value = 3.95
$ws.Cells.Item(13, 17).Value = 1.83

# Row 15
$ws.Cells.Item(15, 17).Value = 1.66
$ws.Cells.Item(15, 18).Value = 1.54
$ws.Cells.Item(15, 19).Value = 2.6
$ws.Cells.Item(15, 24).Value = 28
$ws.Cells.Item(15, 26).Value = 10.5
$ws.Cells.Item(15, 36).Value = 200

# Row 17
$ws.Cells.Item(17, 6).Value = 1.44
$ws.Cells.Item(17, 9).Value = 11
$ws.Cells.Item(17, 10).Value = 4.6
$ws.Cells.Item(17, 11).Value = 5
$ws.Cells.Item(17, 16).Value = 2.04

# Row 19
$ws.Cells.Item(19, 7).Value = 3.2

# Row 20
$ws.Cells.Item(20, 7).Value = 2.28
$ws.Cells.Item(20, 10).Value = 3.25
$ws.Cells.Item(20, 11).Value = 3.45
$ws.Cells.Item(20, 16).Value = 1.71
$ws.Cells.Item(20, 17).Value = 2.28

$wb.Save()